# Il13-Il2rg.xlsx — "update scripts wuth new tpm"
#
# The NATMI LR-pair export was regenerated from a re-run of the pipeline
# with updated TPM values. The "Resolving-Mac" sending-cluster block
# (original rows 10-13) is gone from the new output, and every remaining
# edge's derived-specificity / weight columns were recomputed against the
# new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the four rows for sending cluster "Resolving-Mac" (old rows 10-13). ---
# Deleting the same row index four times in a row removes what used to be
# rows 10, 11, 12 and 13, since each deletion shifts the rows below it up.
$ws.Rows.Item(10).EntireRow.Delete()
$ws.Rows.Item(10).EntireRow.Delete()
$ws.Rows.Item(10).EntireRow.Delete()
$ws.Rows.Item(10).EntireRow.Delete()

# --- Recomputed TPM-derived values for the remaining rows (2-9). ---
$ws.Range("I2").Value = 0.7099439172299504
$ws.Range("J2").Value = 0.7099439172299504
$ws.Range("M2").Value = 43.68636333333333
$ws.Range("N2").Value = 131.05909
$ws.Range("O2").Value = 0.4655970190670373
$ws.Range("P2").Value = 0.4655970190670373
$ws.Range("Q2").Value = 18.59067366801555
$ws.Range("R2").Value = 167.31606301214
$ws.Range("S2").Value = 0.3305477715670403
$ws.Range("T2").Value = 0.3305477715670404

$ws.Range("I3").Value = 0.7099439172299504
$ws.Range("J3").Value = 0.7099439172299504
$ws.Range("O3").Value = 0.01284473362738172
$ws.Range("P3").Value = 0.01284473362738172
$ws.Range("S3").Value = 0.009119040507198648
$ws.Range("T3").Value = 0.009119040507198648

$ws.Range("I4").Value = 0.7099439172299504
$ws.Range("J4").Value = 0.7099439172299504
$ws.Range("M4").Value = 4.981224333333333
$ws.Range("N4").Value = 14.943673
$ws.Range("O4").Value = 0.05308849315764798
$ws.Range("P4").Value = 0.05308849315764799
$ws.Range("Q4").Value = 2.119753373417555
$ws.Range("R4").Value = 19.077780360758
$ws.Range("S4").Value = 0.03768985279217603
$ws.Range("T4").Value = 0.03768985279217603

$ws.Range("I5").Value = 0.7099439172299504
$ws.Range("J5").Value = 0.7099439172299504
$ws.Range("M5").Value = 43.95590833333333
$ws.Range("N5").Value = 131.867725
$ws.Range("O5").Value = 0.4684697541479331
$ws.Range("P5").Value = 0.4684697541479331
$ws.Range("Q5").Value = 18.70537818337222
$ws.Range("R5").Value = 168.34840365035
$ws.Range("S5").Value = 0.3325872523635354
$ws.Range("T5").Value = 0.3325872523635355

$ws.Range("I6").Value = 0.2900560827700495
$ws.Range("J6").Value = 0.2900560827700495
$ws.Range("M6").Value = 43.68636333333333
$ws.Range("N6").Value = 131.05909
$ws.Range("O6").Value = 0.4655970190670373
$ws.Range("P6").Value = 0.4655970190670373
$ws.Range("Q6").Value = 7.595442188223333
$ws.Range("R6").Value = 68.35897969400999
$ws.Range("S6").Value = 0.1350492474999969
$ws.Range("T6").Value = 0.1350492474999969

$ws.Range("I7").Value = 0.2900560827700495
$ws.Range("J7").Value = 0.2900560827700495
$ws.Range("O7").Value = 0.01284473362738172
$ws.Range("P7").Value = 0.01284473362738172
$ws.Range("S7").Value = 0.00372569312018307
$ws.Range("T7").Value = 0.00372569312018307

$ws.Range("I8").Value = 0.2900560827700495
$ws.Range("J8").Value = 0.2900560827700495
$ws.Range("M8").Value = 4.981224333333333
$ws.Range("N8").Value = 14.943673
$ws.Range("O8").Value = 0.05308849315764798
$ws.Range("P8").Value = 0.05308849315764799
$ws.Range("Q8").Value = 0.8660506062663333
$ws.Range("R8").Value = 7.794455456396999
$ws.Range("S8").Value = 0.01539864036547195
$ws.Range("T8").Value = 0.01539864036547195

$ws.Range("I9").Value = 0.2900560827700495
$ws.Range("J9").Value = 0.2900560827700495
$ws.Range("M9").Value = 43.95590833333333
$ws.Range("N9").Value = 131.867725
$ws.Range("O9").Value = 0.4684697541479331
$ws.Range("P9").Value = 0.4684697541479331
$ws.Range("Q9").Value = 7.642306090558333
$ws.Range("R9").Value = 68.780754815025
$ws.Range("S9").Value = 0.1358825017843976
$ws.Range("T9").Value = 0.1358825017843976
